$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain (non-numeric-looking) Price (D) updates - safe to assign directly as text
$ws.Range('D2').Value = '30.383.96'
$ws.Range('D3').Value = '2.004.63'
$ws.Range('D13').Value = '1.999.02'
$ws.Range('D23').Value = '30.435.49'
$ws.Range('D26').Value = '2.228.30'

# Numeric-looking Price (D) updates - force text type via NumberFormat "@" trick
# (applied per-cell, since multi-area ranges only apply to the first area in this host),
# then restore the Normal style so no residual formatting is left on the cell.
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = "Normal"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '324.46'
$ws.Range('D5').Style = "Normal"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5092'
$ws.Range('D7').Style = "Normal"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4139'
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08749'
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.133'
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '43.04'
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '24.55'
$ws.Range('D12').Style = "Normal"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.564'
$ws.Range('D14').Style = "Normal"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.002'
$ws.Range('D16').Style = "Normal"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '94.24'
$ws.Range('D17').Style = "Normal"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06515'
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.90'
$ws.Range('D20').Style = "Normal"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.207'
$ws.Range('D22').Style = "Normal"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.92'
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.224'
$ws.Range('D25').Style = "Normal"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '22.28'
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '162.77'
$ws.Range('D28').Style = "Normal"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.407'
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '131.15'
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.132'
$ws.Range('D31').Style = "Normal"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.068'
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.831'
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.346'
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02520'
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.424'
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06588'
$ws.Range('D38').Style = "Normal"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.2194'
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '9.022'
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.6632'
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.228'
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.56'
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.6156'
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.184'
$ws.Range('D46').Style = "Normal"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '124.17'
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '80.21'
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06886'
$ws.Range('D51').Style = "Normal"

# Volume(1h) (E) updates - these are never numeric-looking (leading/trailing spaces + %)
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('E3').Value = '  +4.82%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('E5').Value = '  +1.45%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E7').Value = '  +1.35%  '
$ws.Range('E8').Value = '  +2.91%  '
$ws.Range('E9').Value = '  +6.38%  '
$ws.Range('E10').Value = '  +2.40%  '
$ws.Range('E11').Value = '  +2.48%  '
$ws.Range('E12').Value = '  +3.56%  '
$ws.Range('E13').Value = '  +4.44%  '
$ws.Range('E14').Value = '  +2.45%  '
$ws.Range('E15').Value = '  +2.34%  '
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('E17').Value = '  +2.33%  '
$ws.Range('E18').Value = '  +1.73%  '
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('E20').Value = '  +4.33%  '
$ws.Range('E22').Value = '  +4.67%  '
$ws.Range('E23').Value = '  +1.27%  '
$ws.Range('E24').Value = '  +5.71%  '
$ws.Range('E25').Value = '  +1.22%  '
$ws.Range('E26').Value = '  +4.45%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('E28').Value = '  +0.73%  '
$ws.Range('E29').Value = '  +6.29%  '
$ws.Range('E30').Value = '  +1.77%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  +1.39%  '
$ws.Range('E33').Value = '  +1.08%  '
$ws.Range('E34').Value = '  +1.64%  '
$ws.Range('E35').Value = '  +12.62%  '
$ws.Range('E36').Value = '  +3.35%  '
$ws.Range('E37').Value = '  +1.68%  '
$ws.Range('E38').Value = '  +2.68%  '
$ws.Range('E39').Value = '  +9.08%  '
$ws.Range('E40').Value = '  +1.56%  '
$ws.Range('E41').Value = '  +1.88%  '
$ws.Range('E42').Value = '  +2.89%  '
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('E44').Value = '  +1.90%  '
$ws.Range('E45').Value = '  +2.67%  '
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('E47').Value = '  +0.88%  '
$ws.Range('E48').Value = '  +4.08%  '
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('E50').Value = '  +1.85%  '
$ws.Range('E51').Value = '  +1.37%  '
